$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.851.82'
$ws.Range('E2').Value = '  +2.15%  '

$ws.Range('D3').Value = '3.728.85'
$ws.Range('E3').Value = '  -1.14%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.64%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.14%  '

$ws.Range('D7').Value = '3.726.75'
$ws.Range('E7').Value = '  -1.03%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.58%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.165'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.02%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.30'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.25%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.15'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.01%  '

$ws.Range('D15').Value = '4.350.36'
$ws.Range('E15').Value = '  -1.11%  '

$ws.Range('D16').Value = '3.730.57'
$ws.Range('E16').Value = '  -1.10%  '

$ws.Range('D17').Value = '68.822.81'
$ws.Range('E17').Value = '  +2.01%  '

$ws.Range('E18').Value = '  +2.35%  '

$ws.Range('E19').Value = '  +0.39%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.19'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +7.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '496.40'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.42'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.722'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.46%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.71'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.92%  '

$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.30'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.25%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000141'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.45%  '

$ws.Range('E27').Value = '  +0.98%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.74%  '

$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('E30').Value = '  +1.13%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.41'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.92'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.83%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.64'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.95%  '

$ws.Range('D34').Value = '3.870.37'
$ws.Range('E34').Value = '  -1.07%  '

$ws.Range('E35').Value = '  +0.81%  '

$ws.Range('D36').Value = '3.663.15'
$ws.Range('E36').Value = '  -1.27%  '

$ws.Range('E37').Value = '  +0.09%  '

$ws.Range('E38').Value = '  +0.74%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.77'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.22%  '

$ws.Range('E40').Value = '  -1.82%  '

$ws.Range('E41').Value = '  +0.35%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '436.03'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.88%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.95'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.36%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.98'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.11%  '

$ws.Range('E45').Value = '  +0.85%  '

$ws.Range('E46').Value = '  +1.80%  '

$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '143.64'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.73%  '

$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '40.45'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0351'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.25%  '

$ws.Range('D51').Value = '2.743.63'
$ws.Range('E51').Value = '  -2.81%  '
